$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append "_123" (destination well count) to the destination-barcode values
# in column E. Rows 2-5 hold "ssdest000000141jul17" and rows 6-7 hold
# "ssdest000000141jul17_384" -- both get "_123" appended.
$ws.Range("E2").Value = "ssdest000000141jul17_123"
$ws.Range("E3").Value = "ssdest000000141jul17_123"
$ws.Range("E4").Value = "ssdest000000141jul17_123"
$ws.Range("E5").Value = "ssdest000000141jul17_123"
$ws.Range("E6").Value = "ssdest000000141jul17_384_123"
$ws.Range("E7").Value = "ssdest000000141jul17_384_123"

# Move the active selection to E7 (was E9).
$ws.Range("E7").Select()
